$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$win = $excel.Windows.Item(1)
$win.ScrollRow = 3
$win.ScrollColumn = 1
Write-Host $win.ScrollRow
